# "Generate Report for Handback"
# Updates the localization-status report after a handback event:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - A "Latest Target File" (F) and "Latest Handback File" (G) hyperlink pair
#    is recorded for every row, on both the zh-cn and de-de sheets
#  - The "Latest Handback DateTime" (H) is stamped with the real handback time

$wb = $excel.ActiveWorkbook

# Excel's blue hyperlink color (RGB 0x6495ED) expressed as the BGR integer
# that the Font.Color COM property expects.
$hyperlinkColor = 15570276

function Set-HandbackRow($ws, $row, $mdName, $mdUrl, $xlfName, $xlfUrl, $handbackDateTime) {
    # Status: handoff is done, the file has been handed back.
    $ws.Cells.Item($row, 3).Value = "Handed back: in sync with en-US"

    # F = Latest Target File (same file that was targeted for translation)
    $fCell = $ws.Cells.Item($row, 6)
    $ws.Hyperlinks.Add($fCell, $mdUrl, "", "", $mdName)
    $fCell.Font.Underline = 2
    $fCell.Font.Color = $hyperlinkColor

    # G = Latest Handback File (the translated file that was handed back)
    $gCell = $ws.Cells.Item($row, 7)
    $ws.Hyperlinks.Add($gCell, $xlfUrl, "", "", $xlfName)
    $gCell.Font.Underline = 2
    $gCell.Font.Color = $hyperlinkColor

    # H = Latest Handback DateTime
    $ws.Cells.Item($row, 8).Value = $handbackDateTime
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Re-create the hyperlinks for the existing columns (A, D) plus the new ones
# (F, G) so that every hyperlinked cell on the sheet is re-added in a single,
# consistent left-to-right / top-to-bottom order.
$wsZh.Hyperlinks.Delete()

Set-HandbackRow $wsZh 2 `
    "844ddcb9-8cde-42e3-bcc3-7fcce877b59e.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/ffbd0d0bbc5695c415791bf0fec43071b6cf69e9/e2e/844ddcb9-8cde-42e3-bcc3-7fcce877b59e.md" `
    "844ddcb9-8cde-42e3-bcc3-7fcce877b59e.294d3a010558aefcb307509f0b8911c3fe6b8321.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8b5fee26ffb7dc509689b0aa3627123140ffa00e/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/hb/844ddcb9-8cde-42e3-bcc3-7fcce877b59e.294d3a010558aefcb307509f0b8911c3fe6b8321.zh-cn.xlf" `
    "2016-03-23 04:11:29"

Set-HandbackRow $wsZh 3 `
    "f91e3bba-44a7-4a39-9ce2-219d3a95d7fa.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/ffbd0d0bbc5695c415791bf0fec43071b6cf69e9/e2e/f91e3bba-44a7-4a39-9ce2-219d3a95d7fa.md" `
    "f91e3bba-44a7-4a39-9ce2-219d3a95d7fa.966dbe3d47744e59447520071ae1fdae09fc56a5.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8b5fee26ffb7dc509689b0aa3627123140ffa00e/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/hb/f91e3bba-44a7-4a39-9ce2-219d3a95d7fa.966dbe3d47744e59447520071ae1fdae09fc56a5.zh-cn.xlf" `
    "2016-03-23 04:11:29"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()

Set-HandbackRow $wsDe 2 `
    "844ddcb9-8cde-42e3-bcc3-7fcce877b59e.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/ffbd0d0bbc5695c415791bf0fec43071b6cf69e9/e2e/844ddcb9-8cde-42e3-bcc3-7fcce877b59e.md" `
    "844ddcb9-8cde-42e3-bcc3-7fcce877b59e.294d3a010558aefcb307509f0b8911c3fe6b8321.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dfd6a4ca7c33f81e97fb71641bdc949dbacd8144/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/hb/844ddcb9-8cde-42e3-bcc3-7fcce877b59e.294d3a010558aefcb307509f0b8911c3fe6b8321.de-de.xlf" `
    "2016-03-23 04:11:44"

Set-HandbackRow $wsDe 3 `
    "f91e3bba-44a7-4a39-9ce2-219d3a95d7fa.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/ffbd0d0bbc5695c415791bf0fec43071b6cf69e9/e2e/f91e3bba-44a7-4a39-9ce2-219d3a95d7fa.md" `
    "f91e3bba-44a7-4a39-9ce2-219d3a95d7fa.966dbe3d47744e59447520071ae1fdae09fc56a5.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dfd6a4ca7c33f81e97fb71641bdc949dbacd8144/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/hb/f91e3bba-44a7-4a39-9ce2-219d3a95d7fa.966dbe3d47744e59447520071ae1fdae09fc56a5.de-de.xlf" `
    "2016-03-23 04:11:44"
